$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '96.872.87'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.53%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.712.76'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +3.40%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.83'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.31%  '
$ws.Range("E6").Value = '  +9.14%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '657.32'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.61%  '
$ws.Range("E8").Value = '  -0.53%  '
$ws.Range("E9").Value = '  +3.55%  '
$ws.Range("E10").Value = '  +0.01%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.707.33'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.31%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '45.46'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.44%  '
$ws.Range("E13").Value = '  +0.98%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.91'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +6.23%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.402.22'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.43%  '
$ws.Range("E16").Value = '  +3.53%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '96.751.90'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.34%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '9.10'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +17.17%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.717.70'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.86%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.22'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.66%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.95'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.58%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.533'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.32%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '526.47'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.28%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.52'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.42%  '
$ws.Range("E25").Value = '  +2.19%  '
$ws.Range("E26").Value = '  -0.74%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '102.74'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.20%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '13.46'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.94%  '
$ws.Range("E29").Value = '  -4.51%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '12.71'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.58%  '
$ws.Range("E31").Value = '  +3.16%  '
$ws.Range("E33").Value = '  +14.50%  '
$ws.Range("E34").Value = '  -0.87%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '32.88'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.66%  '
$ws.Range("E36").Value = '  +0.36%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '658.28'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.69%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.603'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.64%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '9.01'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.76%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.16'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +16.63%  '
$ws.Range("E41").Value = '  +5.57%  '
$ws.Range("E42").Value = '  +3.38%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.974'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.34%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '39.07'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +19.57%  '
$ws.Range("E45").Value = '  +0.04%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0462'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.21%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.445'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.09%  '
$ws.Range("E48").Value = '  +0.41%  '
$ws.Range("E49").Value = '  -0.09%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.72'
$ws.Range("D50").Style = "Normal"
$ws.Range("E51").Value = '  +2.72%  '
